$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 7740
$ws.Range("F6").Value = 97
$ws.Range("F8").Value = 2107
$ws.Range("F9").Value = 8566
$ws.Range("F13").Value = 5717
$ws.Range("F14").Value = 62
$ws.Range("F15").Value = 2669
$ws.Range("F17").Value = 4599
$ws.Range("F18").Value = 355
$ws.Range("F20").Value = 98
$ws.Range("F21").Value = 36
$ws.Range("F22").Value = 557
$ws.Range("F23").Value = 3700
$ws.Range("F24").Value = 72
$ws.Range("F25").Value = 45
$ws.Range("F26").Value = 35
$ws.Range("F27").Value = 10
$ws.Range("F28").Value = 3210
$ws.Range("F30").Value = 290
$ws.Range("F32").Value = 364
$ws.Range("F34").Value = 343
$ws.Range("F35").Value = 1017
$ws.Range("F36").Value = 679
$ws.Range("F37").Value = 12
$ws.Range("F39").Value = 2679
$ws.Range("F43").Value = 3210

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 111
$ws.Range("F3").Value = 137
$ws.Range("F4").Value = 12
$ws.Range("F5").Value = 57
$ws.Range("F7").Value = 40

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1345

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1345
$ws.Range("F5").Value = 7740
$ws.Range("F6").Value = 97
$ws.Range("F8").Value = 2107
$ws.Range("F9").Value = 8566
$ws.Range("F12").Value = 5717
$ws.Range("F13").Value = 62
$ws.Range("F14").Value = 2669
$ws.Range("F16").Value = 4599
$ws.Range("F17").Value = 410
$ws.Range("F18").Value = 98
$ws.Range("F19").Value = 111
$ws.Range("F20").Value = 36
$ws.Range("F21").Value = 137
$ws.Range("F22").Value = 558
$ws.Range("F23").Value = 12
$ws.Range("F24").Value = 3700
$ws.Range("F25").Value = 72
$ws.Range("F26").Value = 45
$ws.Range("F27").Value = 35
$ws.Range("F28").Value = 10
$ws.Range("F29").Value = 3210
$ws.Range("F31").Value = 364
$ws.Range("F33").Value = 343
$ws.Range("F34").Value = 57
$ws.Range("F35").Value = 1017
$ws.Range("F36").Value = 679
$ws.Range("F37").Value = 12
$ws.Range("F39").Value = 40
$ws.Range("F40").Value = 2680
$ws.Range("F44").Value = 3210
